$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix spelling: "Piętka" -> "Piątka"
$ws.Range("A2").Value = "10. Łańcucka Piątka"

# Update selection to B6
$ws.Range("B6").Select()
